# Fixed missing department/region leaders
# Adds a "Department Leader" row (right after the header row) and a
# "Regional Leader" row (at the end of the table) to the attrition table
# on the first worksheet, then grows Table3 to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert "Department Leader" as the new row 2 -------------------------
# (pushes the existing Analyst..Project Manager rows down by one)
$ws.Rows.Item(2).Insert()

# year column ("2009") needs to stay text, like the rest of the column, so
# copy it from a neighboring cell instead of typing it (typing a numeric
# looking string gets auto-converted to a real number by Excel).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial()

$ws.Cells.Item(2, 2).Value = "Department Leader"
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 7

# TerminationRate is stored as text too ("28.6%"). Write it as a formula
# that evaluates to that literal string, then convert the formula to a
# plain value in place so the cell ends up as ordinary text (no percent
# number format / style gets introduced).
$ws.Cells.Item(2, 6).Formula = "=""28.6%"""
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

# --- Append "Regional Leader" as the new last row -------------------------
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial()

$ws.Cells.Item(9, 2).Value = "Regional Leader"
$ws.Cells.Item(9, 3).Value = 14
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 14

# Row 8 (Project Manager) already shows "0.0%" as text; reuse it instead of
# retyping the percentage so the new cell's type/value stay consistent.
$ws.Range("F8").Copy()
$ws.Range("F9").PasteSpecial()

# --- Grow Table3 so the autofilter / table definition covers A1:F9 --------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F9"))

$excel.CutCopyMode = $false
